$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price / volume(1h) data (and two row name/link/price/volume swaps)
$ws.Range("D2").Value = "43.188.87"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "2.372.86"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'303.55"
$ws.Range("D6").Value = "'97.06"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").Value = "'34.21"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E12").Value = "  +3.10%  "
$ws.Range("E13").Value = "  -3.51%  "
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "2.738.69"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "2.357.53"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "43.191.66"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "'12.39"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'6.30"
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("D21").Value = "0.0₃0889"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "'68.28"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "'235.94"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "'31.55"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Value = "'0.0729"
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("D34").Value = "'17.30"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").Value = "'2.30"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'22.88"
$ws.Range("E39").Value = "  +12.70%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.79"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.947.61"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'102.77"
$ws.Range("E43").Value = "  -38.11%  "
$ws.Range("D44").Value = "'0.0281"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  +4.24%  "
$ws.Range("D46").Value = "'9.46"
$ws.Range("E46").Value = "  -9.50%  "
$ws.Range("D47").Value = "'2.76"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").Value = "2.594.87"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").Value = "'72.32"
$ws.Range("E51").Value = "  +1.07%  "
